$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting the existing rows 3-5 down to 4-6.
$ws.Rows.Item(3).Insert()

# Fill in the new row 3 with the new weekly price entry.
$ws.Cells.Item(3, 1).Value = 8
$ws.Cells.Item(3, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(3, 3).Value = "Coquimbo"
$ws.Cells.Item(3, 4).Value = 44452
$ws.Cells.Item(3, 5).Value = 4
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100101
$ws.Cells.Item(3, 8).Value = "Berries"
$ws.Cells.Item(3, 9).Value = 100101001
$ws.Cells.Item(3, 10).Value = "Arándano (blue)"
$ws.Cells.Item(3, 11).Value = "Sin especificar"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 200
$ws.Cells.Item(3, 14).Value = 13000
$ws.Cells.Item(3, 15).Value = 14000
$ws.Cells.Item(3, 16).Value = 13500
$ws.Cells.Item(3, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(3, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(3, 19).Value = 6750
$ws.Cells.Item(3, 20).Value = 2
